# Update rewrite rule numbers in the filled-in table (rows 16-22) to match
# the cell number shown in the index table above, and update the "Further
# explanation" index column (P16:P28) to match those same cell numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 ---
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = 10
$ws.Range("J16").Value = 12
$ws.Range("L16").Value = 14
$ws.Range("M16").Value = 15
$ws.Range("P16").Value = 8

# --- Row 17 ---
$ws.Range("P17").Value = 41

# --- Row 18 ---
$ws.Range("P18").Value = 45

# --- Row 19 ---
$ws.Range("M19").Value = 39
$ws.Range("P19").Value = 10

# --- Row 20 ---
$ws.Range("F20").Value = 40
$ws.Range("G20").Value = 41
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 45
$ws.Range("M20").Value = 47
$ws.Range("P20").Value = 14

# --- Row 21 ---
$ws.Range("P21").Value = 57

# --- Row 22 ---
$ws.Range("F22").Value = 56
$ws.Range("G22").Value = 57
$ws.Range("K22").Value = 61
$ws.Range("O22").Value = 15
$ws.Range("P22").Value = 56

# --- Row 23 ---
$ws.Range("P23").Value = 40

# --- Row 24 ---
$ws.Range("P24").Value = 12

# --- Row 25 ---
$ws.Range("P25").Value = 39

# --- Row 26 ---
$ws.Range("P26").Value = 47

# --- Row 27 ---
$ws.Range("P27").Value = 61

# --- Row 28 (new cell) ---
# Copy the formatting used by the rest of the P column (right-aligned /
# "Footnote number" style, same as P16:P27) onto the newly used P28 cell,
# then set its value.
$ws.Range("P27").Copy($ws.Range("P28"))
$ws.Range("P28").Value = 44

# --- Update selection / active cell to M15 ---
$excel.Goto($ws.Range("M15"), $true)
